$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 76
$ws_ALC.Range("H76").Value = 20004200
$ws_ALC.Range("I76").Value = 25003500
$ws_ALC.Range("K76").Value = 25003500
$ws_ALC.Range("M76").Value = -25003185

# ALC row 79
$ws_ALC.Range("H79").Value = 20004200
$ws_ALC.Range("I79").Value = 25003500
$ws_ALC.Range("K79").Value = 25003500
$ws_ALC.Range("M79").Value = -25002408

# ALC row 138
$ws_ALC.Range("H138").Value = 4091.4849
$ws_ALC.Range("I138").Value = 1005.19446
$ws_ALC.Range("J138").Value = 7795.033
$ws_ALC.Range("K138").Value = 3015.58338
$ws_ALC.Range("L138").Value = 23385.099
$ws_ALC.Range("M138").Value = 2124.41662
$ws_ALC.Range("N138").Value = -33665.099

# ALC row 141
$ws_ALC.Range("H141").Value = 5953903.5
$ws_ALC.Range("I141").Value = 8772723
$ws_ALC.Range("K141").Value = 26318169
$ws_ALC.Range("M141").Value = -26312989

# ARM row 32
$ws_ARM.Range("H32").Value = 1443.01
$ws_ARM.Range("I32").Value = 1452.8247
$ws_ARM.Range("J32").Value = 1125.6666
$ws_ARM.Range("K32").Value = 1452.8247
$ws_ARM.Range("L32").Value = 1125.6666
$ws_ARM.Range("M32").Value = -1165.8247
$ws_ARM.Range("N32").Value = -1699.6666

# ARM row 102
$ws_ARM.Range("H102").Value = 982.2143
$ws_ARM.Range("I102").Value = 986.3
$ws_ARM.Range("K102").Value = 986.3
$ws_ARM.Range("M102").Value = 635.7

# ARM row 105
$ws_ARM.Range("H105").Value = 226999.5
$ws_ARM.Range("J105").Value = 226999.5
$ws_ARM.Range("L105").Value = 226999.5
$ws_ARM.Range("N105").Value = -233987.5

# ARM row 132
$ws_ARM.Range("H132").Value = 9188.034
$ws_ARM.Range("I132").Value = 8625.385
$ws_ARM.Range("J132").Value = 9645.1875
$ws_ARM.Range("K132").Value = 25876.155
$ws_ARM.Range("L132").Value = 28935.5625
$ws_ARM.Range("M132").Value = -23346.155
$ws_ARM.Range("N132").Value = -33995.5625

# CRP row 31
$ws_CRP.Range("H31").Value = 5465.92
$ws_CRP.Range("I31").Value = 2285.7144
$ws_CRP.Range("K31").Value = 2285.7144
$ws_CRP.Range("M31").Value = -1990.7144

# CRP row 34
$ws_CRP.Range("H34").Value = 5465.92
$ws_CRP.Range("I34").Value = 2285.7144
$ws_CRP.Range("K34").Value = 2285.7144
$ws_CRP.Range("M34").Value = -2083.7144

# CRP row 58
$ws_CRP.Range("H58").Value = 8069060
$ws_CRP.Range("I58").Value = 12196738
$ws_CRP.Range("K58").Value = 12196738
$ws_CRP.Range("M58").Value = -12196535

# CRP row 132
$ws_CRP.Range("H132").Value = 4477.0547
$ws_CRP.Range("I132").Value = 2041.0526
$ws_CRP.Range("J132").Value = 9922.235000000001
$ws_CRP.Range("K132").Value = 6123.1578
$ws_CRP.Range("L132").Value = 29766.705
$ws_CRP.Range("M132").Value = -3593.1578
$ws_CRP.Range("N132").Value = -34826.705

# CRP row 136
$ws_CRP.Range("H136").Value = 8069060
$ws_CRP.Range("I136").Value = 12196738
$ws_CRP.Range("K136").Value = 36590214
$ws_CRP.Range("M136").Value = -36587664

# CUL row 18
$ws_CUL.Range("H18").Value = 173.85715
$ws_CUL.Range("I18").Value = 191.33333
$ws_CUL.Range("J18").Value = 69
$ws_CUL.Range("K18").Value = 573.99999
$ws_CUL.Range("L18").Value = 207
$ws_CUL.Range("M18").Value = -404.99999
$ws_CUL.Range("N18").Value = -545

# CUL row 63
$ws_CUL.Range("H63").Value = 500
$ws_CUL.Range("I63").Value = 500
$ws_CUL.Range("K63").Value = 1500
$ws_CUL.Range("M63").Value = -751

# CUL row 66
$ws_CUL.Range("H66").Value = 500
$ws_CUL.Range("I66").Value = 500
$ws_CUL.Range("K66").Value = 4500
$ws_CUL.Range("M66").Value = -756

# CUL row 75
$ws_CUL.Range("H75").Value = 66675336
$ws_CUL.Range("I75").Value = 111115030
$ws_CUL.Range("K75").Value = 333345090
$ws_CUL.Range("M75").Value = -333344092

# CUL row 78
$ws_CUL.Range("H78").Value = 66675336
$ws_CUL.Range("I78").Value = 111115030
$ws_CUL.Range("K78").Value = 1000035270
$ws_CUL.Range("M78").Value = -1000030278

# CUL row 87
$ws_CUL.Range("H87").Value = 1453.25
$ws_CUL.Range("I87").Value = 937.6667
$ws_CUL.Range("K87").Value = 2813.0001
$ws_CUL.Range("M87").Value = -1565.0001

# CUL row 90
$ws_CUL.Range("H90").Value = 1453.25
$ws_CUL.Range("I90").Value = 937.6667
$ws_CUL.Range("K90").Value = 8439.0003
$ws_CUL.Range("M90").Value = -2199.0003

# CUL row 132
$ws_CUL.Range("H132").Value = 12428.429
$ws_CUL.Range("J132").Value = 20666.334
$ws_CUL.Range("L132").Value = 185997.006
$ws_CUL.Range("N132").Value = -191057.006

# CUL row 139
$ws_CUL.Range("H139").Value = 3551.1765
$ws_CUL.Range("I139").Value = 2169.2856
$ws_CUL.Range("K139").Value = 6507.8568
$ws_CUL.Range("M139").Value = -1367.8568

# GSM row 80
$ws_GSM.Range("H80").Value = 2700.75
$ws_GSM.Range("I80").Value = 2802.8333
$ws_GSM.Range("J80").Value = 2394.5
$ws_GSM.Range("K80").Value = 2802.8333
$ws_GSM.Range("L80").Value = 2394.5
$ws_GSM.Range("M80").Value = -1804.8333
$ws_GSM.Range("N80").Value = -4390.5

# GSM row 83
$ws_GSM.Range("H83").Value = 2700.75
$ws_GSM.Range("I83").Value = 2802.8333
$ws_GSM.Range("J83").Value = 2394.5
$ws_GSM.Range("K83").Value = 14014.1665
$ws_GSM.Range("L83").Value = 11972.5
$ws_GSM.Range("M83").Value = -9022.166499999999
$ws_GSM.Range("N83").Value = -21956.5

# GSM row 132
$ws_GSM.Range("H132").Value = 9169.237999999999
$ws_GSM.Range("I132").Value = 3545
$ws_GSM.Range("J132").Value = 12630.308
$ws_GSM.Range("K132").Value = 10635
$ws_GSM.Range("L132").Value = 37890.924
$ws_GSM.Range("M132").Value = -8105
$ws_GSM.Range("N132").Value = -42950.924

# LTW row 46
$ws_LTW.Range("H46").Value = 15875971
$ws_LTW.Range("J46").Value = 18521800
$ws_LTW.Range("L46").Value = 18521800
$ws_LTW.Range("N46").Value = -18522176

# LTW row 68
$ws_LTW.Range("H68").Value = 4294.769
$ws_LTW.Range("I68").Value = 2306.2
$ws_LTW.Range("J68").Value = 5537.625
$ws_LTW.Range("K68").Value = 2306.2
$ws_LTW.Range("L68").Value = 5537.625
$ws_LTW.Range("M68").Value = -1557.2
$ws_LTW.Range("N68").Value = -7035.625

# LTW row 71
$ws_LTW.Range("H71").Value = 4294.769
$ws_LTW.Range("I71").Value = 2306.2
$ws_LTW.Range("J71").Value = 5537.625
$ws_LTW.Range("K71").Value = 11531
$ws_LTW.Range("L71").Value = 27688.125
$ws_LTW.Range("M71").Value = -7787
$ws_LTW.Range("N71").Value = -35176.125

# LTW row 82
$ws_LTW.Range("H82").Value = 1176821.1
$ws_LTW.Range("J82").Value = 2694
$ws_LTW.Range("L82").Value = 2694
$ws_LTW.Range("N82").Value = -3416

# LTW row 85
$ws_LTW.Range("H85").Value = 1176821.1
$ws_LTW.Range("J85").Value = 2694
$ws_LTW.Range("L85").Value = 2694
$ws_LTW.Range("N85").Value = -5190

# LTW row 122
$ws_LTW.Range("H122").Value = 4921.3335
$ws_LTW.Range("I122").Value = 3657.84
$ws_LTW.Range("J122").Value = 6779.4116
$ws_LTW.Range("K122").Value = 10973.52
$ws_LTW.Range("L122").Value = 20338.2348
$ws_LTW.Range("M122").Value = -8523.52
$ws_LTW.Range("N122").Value = -25238.2348

# LTW row 132
$ws_LTW.Range("H132").Value = 10006331
$ws_LTW.Range("I132").Value = 18521356
$ws_LTW.Range("K132").Value = 55564068
$ws_LTW.Range("M132").Value = -55561538

# WVR row 81
$ws_WVR.Range("H81").Value = 66668000
$ws_WVR.Range("I81").Value = 2000
$ws_WVR.Range("J81").Value = 100001000
$ws_WVR.Range("K81").Value = 4000
$ws_WVR.Range("L81").Value = 200002000
$ws_WVR.Range("M81").Value = -2939
$ws_WVR.Range("N81").Value = -200004122

# WVR row 84
$ws_WVR.Range("H84").Value = 66668000
$ws_WVR.Range("I84").Value = 2000
$ws_WVR.Range("J84").Value = 100001000
$ws_WVR.Range("K84").Value = 20000
$ws_WVR.Range("L84").Value = 1000010000
$ws_WVR.Range("M84").Value = -14696
$ws_WVR.Range("N84").Value = -1000020608

# WVR row 122
$ws_WVR.Range("H122").Value = 3232.027
$ws_WVR.Range("I122").Value = 1710.95
$ws_WVR.Range("J122").Value = 5021.5293
$ws_WVR.Range("K122").Value = 5132.85
$ws_WVR.Range("L122").Value = 15064.5879
$ws_WVR.Range("M122").Value = -2682.85
$ws_WVR.Range("N122").Value = -19964.5879

# WVR row 132
$ws_WVR.Range("H132").Value = 14297852
$ws_WVR.Range("I132").Value = 20005570
$ws_WVR.Range("J132").Value = 28555.4
$ws_WVR.Range("K132").Value = 60016710
$ws_WVR.Range("L132").Value = 85666.20000000001
$ws_WVR.Range("M132").Value = -60014180
$ws_WVR.Range("N132").Value = -90726.20000000001

# WVR row 136
$ws_WVR.Range("H136").Value = 25644972
$ws_WVR.Range("I136").Value = 41667560
$ws_WVR.Range("J136").Value = 8828
$ws_WVR.Range("K136").Value = 125002680
$ws_WVR.Range("L136").Value = 26484
$ws_WVR.Range("M136").Value = -125000130
$ws_WVR.Range("N136").Value = -31584
